$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6605
